$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text fix (row 10): server type reference comment
# ---------------------------------------------------------------------------
$ws.Range("H10").Value2 = "server type ref squick/core/base.h"

# ---------------------------------------------------------------------------
# Row 10 height shrink (less wrapped header text now)
# ---------------------------------------------------------------------------
$ws.Rows.Item(10).RowHeight = 30.75

# ---------------------------------------------------------------------------
# Row 11 (MasterServer_1): renumber ports
# ---------------------------------------------------------------------------
$ws.Range("G11").Value2 = 10001
$ws.Range("H11").Value2 = "1"
$ws.Range("J11").Value2 = 10002

# ---------------------------------------------------------------------------
# Row 12 (LoginServer_1): renumber ports / area
# ---------------------------------------------------------------------------
$ws.Range("G12").Value2 = 10010
$ws.Range("H12").Value2 = "2"
$ws.Range("J12").Value2 = 80

# ---------------------------------------------------------------------------
# Row 13 (WorldServer_1): renumber
# ---------------------------------------------------------------------------
$ws.Range("B13").Value2 = "100"
$ws.Range("G13").Value2 = 10101
$ws.Range("H13").Value2 = "5"

# ---------------------------------------------------------------------------
# Row 14 (DbProxyServer_1): renumber
# ---------------------------------------------------------------------------
$ws.Range("B14").Value2 = "300"
$ws.Range("G14").Value2 = 10301
$ws.Range("H14").Value2 = "6"

# ---------------------------------------------------------------------------
# Row 15: used to be GatewayServer_1 -> becomes ProxyServer_1
# (gateway server removed from the table; slot repurposed)
# ---------------------------------------------------------------------------
$ws.Range("A15").Value2 = "ProxyServer_1"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").HorizontalAlignment = -4152
$ws.Range("B15").Value2 = "500"
$ws.Range("C15").Value2 = 5000
$ws.Range("F15").Value2 = "100.100.100.105"
$ws.Range("G15").Value2 = 10501
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").HorizontalAlignment = -4152
$ws.Range("H15").Value2 = "3"
$ws.Range("J15").Value2 = 10502
$ws.Range("K15").Value2 = 10503
$ws.Range("L15").WrapText = $true
$ws.Range("L15").Value2 = 10504

# ---------------------------------------------------------------------------
# Row 16: used to be ProxyServer_1 -> becomes ProxyServer_2
# ---------------------------------------------------------------------------
$ws.Range("A16").Value2 = "ProxyServer_2"
$ws.Range("B16").Value2 = "501"
$ws.Range("G16").Value2 = 10505
$ws.Range("H16").Value2 = "3"
$ws.Range("J16").Value2 = 10506
$ws.Range("K16").Value2 = 10507
$ws.Range("L16").Value2 = 10508

# ---------------------------------------------------------------------------
# Row 17: GameServer_16001 -> GameServer_1
# ---------------------------------------------------------------------------
$ws.Range("A17").Value2 = "GameServer_1"
$ws.Range("B17").Value2 = 1000
$ws.Range("G17").Value2 = 11000
$ws.Range("H17").Value2 = "4"

# ---------------------------------------------------------------------------
# Row 18: GameServer_16002 -> GameServer_2
# ---------------------------------------------------------------------------
$ws.Range("A18").Value2 = "GameServer_2"
$ws.Range("B18").Value2 = 1001
$ws.Range("G18").Value2 = 11001
$ws.Range("H18").Value2 = "4"

# ---------------------------------------------------------------------------
# Row 19: GameplayManagerServer_1 renumber
# ---------------------------------------------------------------------------
$ws.Range("B19").Value2 = 2000
$ws.Range("G19").Value2 = 12000
$ws.Range("H19").Value2 = 7

# ---------------------------------------------------------------------------
# New rows 20 & 21: CdnServer_1 / CdnServer_2
# Copy formats (xlPasteFormats) from the matching row-19 cells so the
# duplicated style indexes line up with the target workbook, then set values.
# ---------------------------------------------------------------------------
$ws.Range("C19").Copy(); $ws.Range("C20").PasteSpecial(-4122)
$ws.Range("D19").Copy(); $ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E19").Copy(); $ws.Range("E20").PasteSpecial(-4122)
$ws.Range("F19").Copy(); $ws.Range("F20").PasteSpecial(-4122)
$ws.Range("I19").Copy(); $ws.Range("I20").PasteSpecial(-4122)
$ws.Range("K19").Copy(); $ws.Range("K20").PasteSpecial(-4122)
$ws.Range("L19").Copy(); $ws.Range("L20").PasteSpecial(-4122)
$ws.Range("M19").Copy(); $ws.Range("M20").PasteSpecial(-4122)
$ws.Range("N19").Copy(); $ws.Range("N20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A20").Value2 = "CdnServer_1"
$ws.Range("B20").Value2 = 3000
$ws.Range("C20").Value2 = 5000
$ws.Range("D20").Value2 = 1
$ws.Range("E20").Value2 = "127.0.0.1"
$ws.Range("F20").Value2 = "100.100.100.108"
$ws.Range("G20").Value2 = 13000
$ws.Range("H20").Value2 = 10
$ws.Range("I20").Value2 = 1
$ws.Range("J20").Value2 = 13001
$ws.Range("K20").Value2 = 0
$ws.Range("L20").Value2 = 0
$ws.Range("M20").Value2 = "sqcuik"
$ws.Range("N20").Value2 = "server_wrold_key"

$ws.Range("C19").Copy(); $ws.Range("C21").PasteSpecial(-4122)
$ws.Range("D19").Copy(); $ws.Range("D21").PasteSpecial(-4122)
$ws.Range("E19").Copy(); $ws.Range("E21").PasteSpecial(-4122)
$ws.Range("F19").Copy(); $ws.Range("F21").PasteSpecial(-4122)
$ws.Range("I19").Copy(); $ws.Range("I21").PasteSpecial(-4122)
$ws.Range("K19").Copy(); $ws.Range("K21").PasteSpecial(-4122)
$ws.Range("L19").Copy(); $ws.Range("L21").PasteSpecial(-4122)
$ws.Range("M19").Copy(); $ws.Range("M21").PasteSpecial(-4122)
$ws.Range("N19").Copy(); $ws.Range("N21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A21").Value2 = "CdnServer_2"
$ws.Range("B21").Value2 = 3001
$ws.Range("C21").Value2 = 5000
$ws.Range("D21").Value2 = 1
$ws.Range("E21").Value2 = "127.0.0.1"
$ws.Range("F21").Value2 = "100.100.100.108"
$ws.Range("G21").Value2 = 13002
$ws.Range("H21").Value2 = 10
$ws.Range("I21").Value2 = 1
$ws.Range("J21").Value2 = 13003
$ws.Range("K21").Value2 = 0
$ws.Range("L21").Value2 = 0
$ws.Range("M21").Value2 = "sqcuik"
$ws.Range("N21").Value2 = "server_wrold_key"

# ---------------------------------------------------------------------------
# Selection moved to B14 (where the login-process key edits were made)
# ---------------------------------------------------------------------------
$ws.Range("B14").Select()
